$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text so values like "1.00" / "0.110"
# keep their exact literal formatting instead of Excel auto-converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.958.88"
$ws.Range("E2").Value = "  -4.60%  "
$ws.Range("D3").Value = "3.075.87"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "538.34"
$ws.Range("E5").Value = "  -6.47%  "
$ws.Range("D6").Value = "131.46"
$ws.Range("E6").Value = "  -12.58%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.069.28"
$ws.Range("E8").Value = "  -4.31%  "
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  -4.64%  "
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  -5.01%  "
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").Value = "  -13.28%  "
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  -5.85%  "
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "34.14"
$ws.Range("E14").Value = "  -10.46%  "
$ws.Range("D15").Value = "3.544.57"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "63.069.46"
$ws.Range("E16").Value = "  -4.64%  "
$ws.Range("D17").Value = "0.110"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "3.078.23"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").Value = "6.55"
$ws.Range("E19").Value = "  -7.68%  "
$ws.Range("D20").Value = "479.93"
$ws.Range("E20").Value = "  -10.10%  "
$ws.Range("D21").Value = "13.20"
$ws.Range("E21").Value = "  -8.98%  "
$ws.Range("D22").Value = "0.697"
$ws.Range("E22").Value = "  -5.78%  "
$ws.Range("D23").Value = "7.13"
$ws.Range("E23").Value = "  -7.29%  "
$ws.Range("D24").Value = "78.06"
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("D25").Value = "11.95"
$ws.Range("E25").Value = "  -11.15%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -8.85%  "
$ws.Range("D28").Value = "8.10"
$ws.Range("E28").Value = "  -13.11%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "25.78"
$ws.Range("E30").Value = "  -6.09%  "
$ws.Range("D31").Value = "1.88"
$ws.Range("E31").Value = "  -16.55%  "
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("D33").Value = "58.30"
$ws.Range("E33").Value = "  +7.05%  "
$ws.Range("D34").Value = "2.39"
$ws.Range("E34").Value = "  -12.40%  "
$ws.Range("D35").Value = "5.90"
$ws.Range("E35").Value = "  -6.72%  "
$ws.Range("D36").Value = "5.18"
$ws.Range("E36").Value = "  -7.28%  "
$ws.Range("D37").Value = "461.60"
$ws.Range("E37").Value = "  -17.65%  "
$ws.Range("D38").Value = "3.103.54"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").Value = "0.0388"
$ws.Range("E39").Value = "  -14.50%  "
$ws.Range("D40").Value = "0.0783"
$ws.Range("E40").Value = "  -8.30%  "
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  -12.06%  "
$ws.Range("D42").Value = "7.99"
$ws.Range("E42").Value = "  -6.46%  "
$ws.Range("D43").Value = "2.49"
$ws.Range("E43").Value = "  -13.71%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.247"
$ws.Range("E45").Value = "  -12.07%  "
$ws.Range("D46").Value = "2.00"
$ws.Range("E46").Value = "  -14.16%  "
$ws.Range("D47").Value = "24.21"
$ws.Range("E47").Value = "  -7.87%  "
$ws.Range("D48").Value = "117.57"
$ws.Range("E48").Value = "  -5.85%  "
$ws.Range("D49").Value = "0.106"
$ws.Range("E49").Value = "  -4.94%  "
$ws.Range("D50").Value = "0.0₃0505"
$ws.Range("E50").Value = "  -8.26%  "
$ws.Range("D51").Value = "1.97"
$ws.Range("E51").Value = "  -10.30%  "

# Restore the default (unstyled) cell style on column D now that the text is committed,
# so no stray formatting is left behind on these data cells.
$ws.Range("D2:D51").Style = "Normal"
